$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.352.55"
$ws.Range("E2").Value = "  +5.04%  "
$ws.Range("D3").Value = "3.628.86"
$ws.Range("E3").Value = "  +9.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.49"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "640.27"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("E7").Value = "  +6.72%  "
$ws.Range("E8").Value = "  +4.81%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  +7.00%  "
$ws.Range("D11").Value = "3.621.05"
$ws.Range("E11").Value = "  +9.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.44"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("E14").Value = "  +5.98%  "
$ws.Range("D15").Value = "4.320.71"
$ws.Range("E15").Value = "  +9.57%  "
$ws.Range("D16").Value = "96.299.71"
$ws.Range("E16").Value = "  +5.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000255"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").Value = "3.631.49"
$ws.Range("E18").Value = "  +9.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  +22.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.25"
$ws.Range("E21").Value = "  +6.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.500"
$ws.Range("E22").Value = "  +11.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "516.27"
$ws.Range("E23").Value = "  +5.23%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000198"
$ws.Range("E25").Value = "  +8.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("E26").Value = "  +10.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "97.41"
$ws.Range("E27").Value = "  +8.61%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.50"
$ws.Range("E28").Value = "  +6.09%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.14"
$ws.Range("E29").Value = "  +21.07%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.58"
$ws.Range("E30").Value = "  +5.42%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.143"
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.182"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.86"
$ws.Range("E35").Value = "  +9.95%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.570"
$ws.Range("E36").Value = "  +8.71%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "574.20"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.86"
$ws.Range("E38").Value = "  +7.66%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  +8.77%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.152"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.930"
$ws.Range("E42").Value = "  +7.70%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0430"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.80"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.71"
$ws.Range("E46").Value = "  +6.14%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.51"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.07"
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.18"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  +4.33%  "
